$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '69.977.10'
$ws.Cells.Item(2, 5).Value = '  -0.14%  '

Set-TextValue 3 4 '3.585.41'
$ws.Cells.Item(3, 5).Value = '  -0.44%  '

$ws.Cells.Item(4, 5).Value = '  +0.13%  '

Set-TextValue 5 4 '579.58'
$ws.Cells.Item(5, 5).Value = '  -1.60%  '

Set-TextValue 6 4 '191.46'
$ws.Cells.Item(6, 5).Value = '  +0.67%  '

$ws.Cells.Item(7, 5).Value = '  -1.80%  '

Set-TextValue 8 4 '3.583.85'
$ws.Cells.Item(8, 5).Value = '  -0.24%  '

Set-TextValue 10 4 '0.181'
$ws.Cells.Item(10, 5).Value = '  +1.86%  '

$ws.Cells.Item(11, 5).Value = '  +0.50%  '

Set-TextValue 12 4 '55.87'
$ws.Cells.Item(12, 5).Value = '  -3.50%  '

Set-TextValue 13 4 '0.0000306'
$ws.Cells.Item(13, 5).Value = '  +5.45%  '

Set-TextValue 14 4 '9.69'
$ws.Cells.Item(14, 5).Value = '  -0.67%  '

Set-TextValue 15 4 '4.165.74'
$ws.Cells.Item(15, 5).Value = '  -0.04%  '

Set-TextValue 16 4 '19.93'
$ws.Cells.Item(16, 5).Value = '  +3.03%  '

Set-TextValue 17 4 '3.587.66'
$ws.Cells.Item(17, 5).Value = '  -0.09%  '

Set-TextValue 18 4 '69.984.45'
$ws.Cells.Item(18, 5).Value = '  +0.08%  '

$ws.Cells.Item(19, 5).Value = '  +1.89%  '

$ws.Cells.Item(20, 5).Value = '  +0.23%  '

$ws.Cells.Item(21, 5).Value = '  -0.31%  '

Set-TextValue 22 4 '475.87'
$ws.Cells.Item(22, 5).Value = '  -3.82%  '

Set-TextValue 23 4 '19.34'
$ws.Cells.Item(23, 5).Value = '  +11.05%  '

Set-TextValue 24 4 '5.03'
$ws.Cells.Item(24, 5).Value = '  -6.52%  '

Set-TextValue 25 4 '4.39'
$ws.Cells.Item(25, 5).Value = '  -1.62%  '

Set-TextValue 26 4 '95.66'
$ws.Cells.Item(26, 5).Value = '  +5.48%  '

Set-TextValue 27 4 '3.00'
$ws.Cells.Item(27, 5).Value = '  -3.55%  '

Set-TextValue 28 4 '11.09'
$ws.Cells.Item(28, 5).Value = '  +0.09%  '

Set-TextValue 29 4 '9.42'
$ws.Cells.Item(29, 5).Value = '  +0.09%  '

Set-TextValue 30 4 '32.28'
$ws.Cells.Item(30, 5).Value = '  +0.11%  '

Set-TextValue 31 4 '7.67'
$ws.Cells.Item(31, 5).Value = '  +0.60%  '

$ws.Cells.Item(32, 5).Value = '  +0.16%  '

$ws.Cells.Item(33, 5).Value = '  +1.45%  '

Set-TextValue 34 4 '66.51'
$ws.Cells.Item(34, 5).Value = '  +2.27%  '

Set-TextValue 35 4 '589.59'
$ws.Cells.Item(35, 5).Value = '  -4.94%  '

Set-TextValue 36 4 '39.01'
$ws.Cells.Item(36, 5).Value = '  +2.32%  '

$ws.Cells.Item(37, 5).Value = '  +0.08%  '

Set-TextValue 38 4 '0.0₃0802'
$ws.Cells.Item(38, 5).Value = '  -1.88%  '

Set-TextValue 39 4 '0.396'
$ws.Cells.Item(39, 5).Value = '  -2.03%  '

Set-TextValue 40 4 '3.24'
$ws.Cells.Item(40, 5).Value = '  +20.01%  '

$ws.Cells.Item(41, 5).Value = '  -5.96%  '

Set-TextValue 42 4 '3.46'
$ws.Cells.Item(42, 5).Value = '  -4.86%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 43 4 '3.234.60'
$ws.Cells.Item(43, 5).Value = '  -2.15%  '

$ws.Cells.Item(44, 2).Value = 'Fetch.AI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 44 4 '2.85'
$ws.Cells.Item(44, 5).Value = '  +6.80%  '

$ws.Cells.Item(45, 5).Value = '  +0.97%  '

Set-TextValue 46 4 '0.0444'
$ws.Cells.Item(46, 5).Value = '  -0.46%  '

Set-TextValue 47 4 '3.35'
$ws.Cells.Item(47, 5).Value = '  +2.54%  '

Set-TextValue 48 4 '9.45'
$ws.Cells.Item(48, 5).Value = '  +3.73%  '

$ws.Cells.Item(49, 5).Value = '  +0.63%  '

$ws.Cells.Item(50, 5).Value = '  +0.18%  '

Set-TextValue 51 4 '3.14'
$ws.Cells.Item(51, 5).Value = '  -5.17%  '
